$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = "hi"
$ws.Range("D8").Select() | Out-Null
